$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cell: D1 = "Total Heures"
$ws.Range("D1").Value = "Total Heures"

# New formula cell: D2 = SUM of the (whole, now 1-row-larger) hours column
$ws.Range("D2").Formula = "=SUM(B2:B13)"

# New data row 13 - entered the same way as the existing rows (text, the
# sheet already suppresses the "number stored as text" warning for this),
# so force text storage with a leading apostrophe like a real user typing
# into the grid would.
$ws.Range("A13").Value = "'01/02/2024"
$ws.Range("B13").Value = "'7"

# New trailing blank row 14 - an empty-but-present cell at B14 (matches the
# widened dimension/ignoredErrors range down to row 14).
$ws.Range("B14").Value = "'"
